$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All cells in columns B-E of this sheet hold text (coin name / URL / price /
# volume%) in the source data, even when a price string happens to look like a
# plain number (e.g. "1.003"). Force text storage via NumberFormat "@" so Excel
# does not silently convert these into numeric cells, then restore the default
# "Normal" style so no stray number-format style lingers on the cell.

$ws.Range("D2").Value = "22.381.98"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "1.567.61"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.004"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3742"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07535"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.130"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.69%  "

$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.928"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.881"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("D16").Value = "1.565.06"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001122"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06733"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.171"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.12%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("D24").Value = "22.370.64"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.380"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.706"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.032"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").Value = "1.741.10"
$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.019"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9845"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.50%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.012"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.423"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08466"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02476"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2276"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06434"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.382"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6262"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.22%  "

$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.801"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5877"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.057"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.80%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.258"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("E51").Value = "  +0.83%  "
